$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (after the header row, before
# the existing row 2), shifting the rest of the data down by two rows.
$ws.Rows("2:3").Insert(-4161)

$ws.Range("A2").Value = "'2026-01-16"
$ws.Range("B2").Value = "Por qué tantas personas de 30 siguen actuando como adolescentes al elegir pareja, según expertos en comportamiento"
$ws.Range("C2").Value = "Diario ADN"
$ws.Range("D2").Value = "Sin identificar"
$ws.Range("E2").Value = "https://www.diarioadn.co/seccion/vida"
$ws.Range("F2").Value = "Afirman que la madurez emocional no llega sola con los años sino que se construye."

$ws.Range("A3").Value = "'2026-01-16"
$ws.Range("B3").Value = 'Yumbo se pone la camiseta del futuro: lanza Aulas STEAM con inversión histórica de $ 9.500 millones'
$ws.Range("C3").Value = "Diario ADN"
$ws.Range("D3").Value = "Sin identificar"
$ws.Range("E3").Value = "https://www.diarioadn.co/seccion/regiones"
$ws.Range("F3").Value = "Iniciativa busca impactar a más de 10.400 estudiantes"

# Drop the formatting the insert carried down from the header row so the
# new rows match the unstyled look of the rest of the data rows.
$ws.Rows("2:3").ClearFormats()
